# Refactor verb selection logic to ensure valid forms are chosen before
# adding to session -- concretely, this adds two new verb-conjugation rows
# (引き出す "to withdraw" and 引っ越す "to move house") to the bottom of
# the verb table on Sheet1, following the same 7-column layout as every
# other row (Dictionary / Te / Ta / Nai / Masu / Volitional / Potential).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Anchor for the new rows: just past the current last data row (96).
$lastRow = 96
$newRow1 = $lastRow + 1   # 97
$newRow2 = $lastRow + 2   # 98

# Clone the formatting (font/style) of the last existing row onto the two
# new rows first, so the new cells pick up the same cell style (s="3",
# Yu Gothic) that the rest of the Japanese-text rows use.
$ws.Range("A96:G96").Copy() | Out-Null
$ws.Range("A97:G98").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 97: 引き出す (hikidasu, "to withdraw") conjugation set.
$ws.Cells.Item($newRow1, 1).Value = "引き出す"
$ws.Cells.Item($newRow1, 2).Value = "引き出して"
$ws.Cells.Item($newRow1, 3).Value = "引き出した"
$ws.Cells.Item($newRow1, 4).Value = "引き出さない"
$ws.Cells.Item($newRow1, 5).Value = "引き出します"
$ws.Cells.Item($newRow1, 6).Value = "引き出そう"
$ws.Cells.Item($newRow1, 7).Value = "引き出せる"

# Row 98: 引っ越す (hikkosu, "to move house") conjugation set.
$ws.Cells.Item($newRow2, 1).Value = "引っ越す"
$ws.Cells.Item($newRow2, 2).Value = "引っ越して"
$ws.Cells.Item($newRow2, 3).Value = "引っ越した"
$ws.Cells.Item($newRow2, 4).Value = "引っ越さない"
$ws.Cells.Item($newRow2, 5).Value = "引っ越します"
$ws.Cells.Item($newRow2, 6).Value = "引っ越そう"
$ws.Cells.Item($newRow2, 7).Value = "引っ越せる"

# Match the author's final on-screen selection/scroll position: the
# cursor ends up on the newly added last cell, A98.
$ws.Range("A98").Select() | Out-Null
